$wb = $excel.ActiveWorkbook

$wsEmail = $wb.Worksheets.Item("Email Template")
$wsValidations = $wb.Worksheets.Item("Field Validations")

# Update the sample e-mail / first-name / last-name data on the Email Template sheet.
$wsEmail.Range("A2").Value = "wasimakramb325@gmail.com"
$wsEmail.Range("B2").Value = "wasim"
$wsEmail.Range("C2").Value = "akram"

$wsEmail.Range("A3").Value = "wakram@dacgroup.com"
$wsEmail.Range("B3").Value = "B"
$wsEmail.Range("C3").Value = "Akram"
$wsEmail.Range("D3").Value = 9990038502

# Update the selected cell on each sheet (Field Validations first, so the
# Email Template sheet ends up as the active/selected tab, matching the
# original workbook).
$wsValidations.Range("B4").Select()
$wsEmail.Range("D3").Select()
